$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.827.88"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.083.68"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.60"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.65"
$ws.Range("E6").Value = "  +4.27%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.081.50"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.81"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.599.20"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.825.23"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.085.56"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.50"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.78"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.28"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.89"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0923"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.944"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.69"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.311"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "48.90"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("E44").Value = "  +8.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.797.40"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "368.18"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.78"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0342"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.66"
$ws.Range("E50").Value = "  +4.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  +6.26%  "
